$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4000  # H62 was 4500
$ws.Cells.Item(62, 9).Value = 4000  # I62 was 4500
$ws.Cells.Item(62, 11).Value = 4000  # K62 was 4500
$ws.Cells.Item(62, 13).Value = -3376  # M62 was -3876
$ws.Cells.Item(65, 8).Value = 4000  # H65 was 4500
$ws.Cells.Item(65, 9).Value = 4000  # I65 was 4500
$ws.Cells.Item(65, 11).Value = 20000  # K65 was 22500
$ws.Cells.Item(65, 13).Value = -16880  # M65 was -19380
$ws.Cells.Item(137, 8).Value = 3705708.5  # H137 was 2704357.5
$ws.Cells.Item(137, 9).Value = 4349582  # I137 was 3126459
$ws.Cells.Item(137, 10).Value = 3435.75  # J137 was 2908.6
$ws.Cells.Item(137, 11).Value = 13048746  # K137 was 9379377
$ws.Cells.Item(137, 12).Value = 10307.25  # L137 was 8725.799999999999
$ws.Cells.Item(137, 13).Value = -13046196  # M137 was -9376827
$ws.Cells.Item(137, 14).Value = -15407.25  # N137 was -13825.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 47715628  # H61 was 34552900
$ws.Cells.Item(61, 9).Value = 62563324  # I61 was 47667476
$ws.Cells.Item(61, 10).Value = 203000  # J61 was 127137.5
$ws.Cells.Item(61, 11).Value = 62563324  # K61 was 47667476
$ws.Cells.Item(61, 12).Value = 203000  # L61 was 127137.5
$ws.Cells.Item(61, 13).Value = -62563112  # M61 was -47667264
$ws.Cells.Item(61, 14).Value = -203424  # N61 was -127561.5
$ws.Cells.Item(74, 8).Value = 10082464  # H74 was 6633710
$ws.Cells.Item(74, 9).Value = 15689381  # I74 was 8966014
$ws.Cells.Item(74, 10).Value = 114611.11  # J74 was 103260
$ws.Cells.Item(74, 11).Value = 15689381  # K74 was 8966014
$ws.Cells.Item(74, 12).Value = 114611.11  # L74 was 103260
$ws.Cells.Item(74, 13).Value = -15688507  # M74 was -8965140
$ws.Cells.Item(74, 14).Value = -116359.11  # N74 was -105008
$ws.Cells.Item(77, 8).Value = 10082464  # H77 was 6633710
$ws.Cells.Item(77, 9).Value = 15689381  # I77 was 8966014
$ws.Cells.Item(77, 10).Value = 114611.11  # J77 was 103260
$ws.Cells.Item(77, 11).Value = 78446905  # K77 was 44830070
$ws.Cells.Item(77, 12).Value = 573055.55  # L77 was 516300
$ws.Cells.Item(77, 13).Value = -78442537  # M77 was -44825702
$ws.Cells.Item(77, 14).Value = -581791.55  # N77 was -525036
$ws.Cells.Item(110, 8).Value = 1001690.4  # H110 was 2002601.4
$ws.Cells.Item(110, 9).Value = 1429470.1  # I110 was 5000505.5
$ws.Cells.Item(110, 10).Value = 3537.6667  # J110 was 3998.6667
$ws.Cells.Item(110, 11).Value = 1429470.1  # K110 was 5000505.5
$ws.Cells.Item(110, 12).Value = 3537.6667  # L110 was 3998.6667
$ws.Cells.Item(110, 13).Value = -1427425.1  # M110 was -4998460.5
$ws.Cells.Item(110, 14).Value = -7627.6667  # N110 was -8088.6667
$ws.Cells.Item(122, 8).Value = 3473888.5  # H122 was 5850117
$ws.Cells.Item(122, 9).Value = 1656.32  # I122 was 2194.625
$ws.Cells.Item(122, 10).Value = 15874718  # J122 was 37039036
$ws.Cells.Item(122, 11).Value = 4968.96  # K122 was 6583.875
$ws.Cells.Item(122, 12).Value = 47624154  # L122 was 111117108
$ws.Cells.Item(122, 13).Value = -2518.96  # M122 was -4133.875
$ws.Cells.Item(122, 14).Value = -47629054  # N122 was -111122008
$ws.Cells.Item(132, 8).Value = 35027.9  # H132 was 30438.25
$ws.Cells.Item(132, 9).Value = 25162.342  # I132 was 20290.04
$ws.Cells.Item(132, 10).Value = 57499.445  # J132 was 60882.883
$ws.Cells.Item(132, 11).Value = 75487.026  # K132 was 60870.12
$ws.Cells.Item(132, 12).Value = 172498.335  # L132 was 182648.649
$ws.Cells.Item(132, 13).Value = -72957.026  # M132 was -58340.12
$ws.Cells.Item(132, 14).Value = -177558.335  # N132 was -187708.649
$ws.Cells.Item(136, 8).Value = 47715628  # H136 was 34552900
$ws.Cells.Item(136, 9).Value = 62563324  # I136 was 47667476
$ws.Cells.Item(136, 10).Value = 203000  # J136 was 127137.5
$ws.Cells.Item(136, 11).Value = 187689972  # K136 was 143002428
$ws.Cells.Item(136, 12).Value = 609000  # L136 was 381412.5
$ws.Cells.Item(136, 13).Value = -187687422  # M136 was -142999878
$ws.Cells.Item(136, 14).Value = -614100  # N136 was -386512.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 45456296  # H105 was 50001748
$ws.Cells.Item(105, 9).Value = 100001400  # I105 was 125001304
$ws.Cells.Item(105, 11).Value = 100001400  # K105 was 125001304
$ws.Cells.Item(105, 13).Value = -99999653  # M105 was -124999557
$ws.Cells.Item(126, 8).Value = 50000  # H126 was 0
$ws.Cells.Item(126, 10).Value = 50000  # J126 was 0
$ws.Cells.Item(126, 12).Value = 50000  # L126 was 0
$ws.Cells.Item(126, 14).Value = -59880  # N126 was None
$ws.Cells.Item(134, 8).Value = 2244.8474  # H134 was 2404.8
$ws.Cells.Item(134, 9).Value = 1761.8541  # I134 was 1850.3778
$ws.Cells.Item(134, 10).Value = 4352.4546  # J134 was 4899.7
$ws.Cells.Item(134, 11).Value = 5285.5623  # K134 was 5551.1334
$ws.Cells.Item(134, 12).Value = 13057.3638  # L134 was 14699.1
$ws.Cells.Item(134, 13).Value = -2750.5623  # M134 was -3016.1334
$ws.Cells.Item(134, 14).Value = -18127.3638  # N134 was -19769.1
$ws.Cells.Item(141, 8).Value = 43106  # H141 was 36274.5
$ws.Cells.Item(141, 9).Value = 40709  # I141 was 32399.334
$ws.Cells.Item(141, 11).Value = 40709  # K141 was 32399.334
$ws.Cells.Item(141, 13).Value = -35529  # M141 was -27219.334
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 22561.75  # H132 was 19655.564
$ws.Cells.Item(132, 9).Value = 1472.8684  # I132 was 1184.2391
$ws.Cells.Item(132, 10).Value = 102699.5  # J132 was 114064.555
$ws.Cells.Item(132, 11).Value = 4418.6052  # K132 was 3552.7173
$ws.Cells.Item(132, 12).Value = 308098.5  # L132 was 342193.665
$ws.Cells.Item(132, 13).Value = -1888.6052  # M132 was -1022.7173
$ws.Cells.Item(132, 14).Value = -313158.5  # N132 was -347253.665
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 735.64  # H5 was 745.4583
$ws.Cells.Item(5, 9).Value = 318.1  # I5 was 297.8889
$ws.Cells.Item(5, 11).Value = 954.3000000000001  # K5 was 893.6667
$ws.Cells.Item(5, 13).Value = -842.3000000000001  # M5 was -781.6667
$ws.Cells.Item(59, 8).Value = 2121.2  # H59 was 3000
$ws.Cells.Item(59, 9).Value = 803  # I59 was 0
$ws.Cells.Item(59, 11).Value = 2409  # K59 was 0
$ws.Cells.Item(59, 13).Value = -1869  # M59 was None
$ws.Cells.Item(132, 8).Value = 954.6429000000001  # H132 was 895
$ws.Cells.Item(132, 9).Value = 686.5  # I132 was 670.38464
$ws.Cells.Item(132, 11).Value = 6178.5  # K132 was 6033.46176
$ws.Cells.Item(132, 13).Value = -3648.5  # M132 was -3503.46176
$ws.Cells.Item(135, 8).Value = 735.64  # H135 was 745.4583
$ws.Cells.Item(135, 9).Value = 318.1  # I135 was 297.8889
$ws.Cells.Item(135, 11).Value = 2862.9  # K135 was 2681.0001
$ws.Cells.Item(135, 13).Value = -327.9000000000001  # M135 was -146.0000999999997
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3703.1  # H80 was 3828.7896
$ws.Cells.Item(80, 9).Value = 2798.3333  # I80 was 3226.6667
$ws.Cells.Item(80, 10).Value = 3862.7646  # J80 was 3941.6875
$ws.Cells.Item(80, 11).Value = 2798.3333  # K80 was 3226.6667
$ws.Cells.Item(80, 12).Value = 3862.7646  # L80 was 3941.6875
$ws.Cells.Item(80, 13).Value = -1800.3333  # M80 was -2228.6667
$ws.Cells.Item(80, 14).Value = -5858.7646  # N80 was -5937.6875
$ws.Cells.Item(83, 8).Value = 3703.1  # H83 was 3828.7896
$ws.Cells.Item(83, 9).Value = 2798.3333  # I83 was 3226.6667
$ws.Cells.Item(83, 10).Value = 3862.7646  # J83 was 3941.6875
$ws.Cells.Item(83, 11).Value = 13991.6665  # K83 was 16133.3335
$ws.Cells.Item(83, 12).Value = 19313.823  # L83 was 19708.4375
$ws.Cells.Item(83, 13).Value = -8999.666499999999  # M83 was -11141.3335
$ws.Cells.Item(83, 14).Value = -29297.823  # N83 was -29692.4375
$ws.Cells.Item(102, 8).Value = 985.3333  # H102 was 800
$ws.Cells.Item(102, 9).Value = 837.3333  # I102 was 800
$ws.Cells.Item(102, 10).Value = 1133.3334  # J102 was 0
$ws.Cells.Item(102, 11).Value = 837.3333  # K102 was 800
$ws.Cells.Item(102, 12).Value = 1133.3334  # L102 was 0
$ws.Cells.Item(102, 13).Value = 784.6667  # M102 was 822
$ws.Cells.Item(102, 14).Value = -4377.3334  # N102 was None
$ws.Cells.Item(132, 8).Value = 44015.297  # H132 was 39861.598
$ws.Cells.Item(132, 9).Value = 31560.395  # I132 was 28223.459
$ws.Cells.Item(132, 10).Value = 73373.28999999999  # J132 was 68569
$ws.Cells.Item(132, 11).Value = 94681.185  # K132 was 84670.37699999999
$ws.Cells.Item(132, 12).Value = 220119.87  # L132 was 205707
$ws.Cells.Item(132, 13).Value = -92151.185  # M132 was -82140.37699999999
$ws.Cells.Item(132, 14).Value = -225179.87  # N132 was -210767
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 722.4545000000001  # H22 was 661.9474
$ws.Cells.Item(22, 9).Value = 461.25  # I22 was 525.9091
$ws.Cells.Item(22, 10).Value = 871.7143  # J22 was 849
$ws.Cells.Item(22, 11).Value = 461.25  # K22 was 525.9091
$ws.Cells.Item(22, 12).Value = 871.7143  # L22 was 849
$ws.Cells.Item(22, 13).Value = -166.25  # M22 was -230.9091
$ws.Cells.Item(22, 14).Value = -1461.7143  # N22 was -1439
$ws.Cells.Item(27, 8).Value = 722.4545000000001  # H27 was 661.9474
$ws.Cells.Item(27, 9).Value = 461.25  # I27 was 525.9091
$ws.Cells.Item(27, 10).Value = 871.7143  # J27 was 849
$ws.Cells.Item(27, 11).Value = 461.25  # K27 was 525.9091
$ws.Cells.Item(27, 12).Value = 871.7143  # L27 was 849
$ws.Cells.Item(27, 13).Value = -354.25  # M27 was -418.9091
$ws.Cells.Item(27, 14).Value = -1085.7143  # N27 was -1063
$ws.Cells.Item(40, 8).Value = 2600.6667  # H40 was 2890.4
$ws.Cells.Item(40, 9).Value = 2600.6667  # I40 was 2888
$ws.Cells.Item(40, 10).Value = 0  # J40 was 2900
$ws.Cells.Item(40, 11).Value = 2600.6667  # K40 was 2888
$ws.Cells.Item(40, 12).Value = 0  # L40 was 2900
$ws.Cells.Item(40, 13).ClearContents()  # M40 was -2752
$ws.Cells.Item(40, 14).Value = -2464.6667  # N40 was -3172
$ws.Cells.Item(46, 8).Value = 869.8570999999999  # H46 was 895
$ws.Cells.Item(46, 9).Value = 772.25  # I46 was 860
$ws.Cells.Item(46, 11).Value = 772.25  # K46 was 860
$ws.Cells.Item(46, 13).Value = -584.25  # M46 was -672
$ws.Cells.Item(55, 8).Value = 348.16666  # H55 was 240.57143
$ws.Cells.Item(55, 9).Value = 132.5  # I55 was 136
$ws.Cells.Item(55, 10).Value = 456  # J55 was 502
$ws.Cells.Item(55, 11).Value = 132.5  # K55 was 136
$ws.Cells.Item(55, 12).Value = 456  # L55 was 502
$ws.Cells.Item(55, 13).Value = 40.5  # M55 was 37
$ws.Cells.Item(55, 14).Value = -802  # N55 was -848
$ws.Cells.Item(122, 8).Value = 3594.077  # H122 was 4479.75
$ws.Cells.Item(122, 9).Value = 2959.6667  # I122 was 3954
$ws.Cells.Item(122, 10).Value = 4137.857  # J122 was 4655
$ws.Cells.Item(122, 11).Value = 8879.000100000001  # K122 was 11862
$ws.Cells.Item(122, 12).Value = 12413.571  # L122 was 13965
$ws.Cells.Item(122, 13).Value = -6429.000100000001  # M122 was -9412
$ws.Cells.Item(122, 14).Value = -17313.571  # N122 was -18865
$ws.Cells.Item(136, 8).Value = 80246.53999999999  # H136 was 49934.69
$ws.Cells.Item(136, 9).Value = 44974.348  # I136 was 27462.025
$ws.Cells.Item(136, 10).Value = 350666.66  # J136 was 263425
$ws.Cells.Item(136, 11).Value = 134923.044  # K136 was 82386.07500000001
$ws.Cells.Item(136, 12).Value = 1051999.98  # L136 was 790275
$ws.Cells.Item(136, 13).Value = -132373.044  # M136 was -79836.07500000001
$ws.Cells.Item(136, 14).Value = -1057099.98  # N136 was -795375
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2864.0557  # H122 was 2172.9092
$ws.Cells.Item(122, 9).Value = 2175.4  # I122 was 1632.6
$ws.Cells.Item(122, 10).Value = 3724.875  # J122 was 3004.1538
$ws.Cells.Item(122, 11).Value = 6526.200000000001  # K122 was 4897.799999999999
$ws.Cells.Item(122, 12).Value = 11174.625  # L122 was 9012.4614
$ws.Cells.Item(122, 13).Value = -4076.200000000001  # M122 was -2447.799999999999
$ws.Cells.Item(122, 14).Value = -16074.625  # N122 was -13912.4614
$ws.Cells.Item(132, 8).Value = 58879.17  # H132 was 51499.074
$ws.Cells.Item(132, 9).Value = 60007.824  # I132 was 40987
$ws.Cells.Item(132, 10).Value = 57813.223  # J132 was 69019.2
$ws.Cells.Item(132, 11).Value = 180023.472  # K132 was 122961
$ws.Cells.Item(132, 12).Value = 173439.669  # L132 was 207057.6
$ws.Cells.Item(132, 13).Value = -177493.472  # M132 was -120431
$ws.Cells.Item(132, 14).Value = -178499.669  # N132 was -212117.6
$ws.Cells.Item(136, 8).Value = 43901.426  # H136 was 33052.793
$ws.Cells.Item(136, 9).Value = 23803.795  # I136 was 19830.584
$ws.Cells.Item(136, 10).Value = 338666.66  # J136 was 103130.5
$ws.Cells.Item(136, 11).Value = 71411.38499999999  # K136 was 59491.75199999999
$ws.Cells.Item(136, 12).Value = 1015999.98  # L136 was 309391.5
$ws.Cells.Item(136, 13).Value = -68861.38499999999  # M136 was -56941.75199999999
$ws.Cells.Item(136, 14).Value = -1021099.98  # N136 was -314491.5
